$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Part 1: "It references NuGet package..." paragraph.
#
#   "  You'll need to install this.  It is used in two projects:"
#     -> "  You'll need to install this" + " as it is" + " is used in "
#        + "multiple projects."
#
#   and fold the two bullet sub-items ("ApiLayer" / "DataAccessLayer")
#   away (their wording is superseded by the new sentence).
#
# NOTE: Find/Replace's ReplaceWith re-types the whole matched span, and
# this runtime's autocorrect turns straight apostrophes into curly ones
# when text is retyped that way. "You'll" must survive untouched, so we
# only ever Find/Replace the apostrophe-free tail of that sentence and
# grow the rest with InsertAfter (confirmed not to trigger autocorrect).
# ---------------------------------------------------------------------

$rng = $d.Content
$rng.Find.Execute(".  It is used in two projects:", $true, $false, $false, $false, $false, $true, 1, $false, "", 2)

$p6 = $d.Paragraphs(6)
$p6.Range.InsertAfter(" as it is")
$p6 = $d.Paragraphs(6)
$p6.Range.InsertAfter(" is used in ")
$p6 = $d.Paragraphs(6)
$p6.Range.InsertAfter("multiple projects.")

# Remove the "ApiLayer" and "DataAccessLayer" bullet paragraphs outright
# (same paragraph index, 7, each time since the prior one disappears).
$d.Paragraphs(7).Range.Delete()
$d.Paragraphs(7).Range.Delete()

# ---------------------------------------------------------------------
# Part 2: collapse the "complete solution" / "To build" / "Using a DOS
# prompt" bullets into a single paragraph. The surviving paragraph keeps
# its original bookmark (_GoBack) sitting between "...Change the JSON
# co" and "ntent as needed...".
#
# NOTE: this runtime's Range.Find.Execute does not stay confined to the
# Range it is called on (it can jump back to an earlier match in the
# story), so once we need a Find restricted to a sub-range we instead
# compute offsets by hand with $d.Range(start, end).
# ---------------------------------------------------------------------

# Insert the first chunk of the "Using a DOS prompt..." text right
# before the _GoBack bookmark (InsertBefore does not trigger the
# autocorrect retype, so "I've" stays a straight apostrophe).
$bm = $d.Bookmarks("_GoBack")
$bm.Range.InsertBefore("Using a DOS prompt invoke the application with arguments as documented below.  You can do this from either the bin/Debug or bin/Release folder.  I've included a JSON folder with the content needed to perform a run.  Change the JSON co")

# Swap out "The complete solution is included as a zip file." (now
# immediately after the bookmark) for the remainder of the sentence.
$rng = $d.Content
$rng.Find.Execute("The complete solution is included as a zip file.")
$rng.Text = "ntent as needed to test out different combinations of Condition/Topography/Date/Room."

# The paragraph that now reads "...Change the JSON co" + bookmark +
# "ntent as needed...Room." is the surviving one; find it via the
# bookmark so later paragraph-index shifts don't matter.
$survivor = $d.Bookmarks("_GoBack").Range.Paragraphs(1)

# Delete the (now redundant) "To build, unzip the solution and perform
# a build." paragraph entirely -- it immediately follows the survivor.
$survivor.Next().Range.Delete()

# The paragraph following the survivor is now the original "Using a DOS
# prompt ... Room." + tab paragraph. Its sentence is now duplicated
# (already merged into the survivor above), so delete just that
# sentence (computed by length, since a sub-range Find is unreliable
# here), keeping the tab run intact.
$survivor = $d.Bookmarks("_GoBack").Range.Paragraphs(1)
$dupPara = $survivor.Next()
$dupSentence = "Using a DOS prompt invoke the application with arguments as documented below.  You can do this from either the bin/Debug or bin/Release folder.  I've included a JSON folder with the content needed to perform a run.  Change the JSON content as needed to test out different combinations of Condition/Topography/Date/Room."
$delStart = $dupPara.Range.Start
$delEnd = $delStart + $dupSentence.Length
$d.Range($delStart, $delEnd).Delete()

# Merge the now-empty-but-for-a-tab paragraph back into the survivor by
# deleting the survivor's own trailing paragraph mark.
$survivor = $d.Bookmarks("_GoBack").Range.Paragraphs(1)
$markStart = $survivor.Range.End - 1
$d.Range($markStart, $survivor.Range.End).Delete()
